# Replace the hard-coded parent/organization names in the report header
# with template placeholders, so the report generator can fill them in
# per-run (ParentOrganizationName / OrganizationName).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 used to hold the hard-coded "ГЛАВНОЕ УПРАВЛЕНИЕ ..." parent org name.
$ws.Range("A1").Value = "{ParentOrganizationName}"

# Row 2 used to hold the hard-coded "УО «Пинский ...»" organization name.
$ws.Range("A2").Value = "{OrganizationName}"

# Leave the selection where the author last left it while editing.
$ws.Range("J7").Select()
